# Add three new vowel rows ("long + pharyngealized" variants) to the
# phonological feature table on Sheet1: iʶː (before "iː"), uʶː (before "uː"),
# and aʶː (before "aː"). Each new row is inserted above its matching
# "long" vowel row so the table keeps its existing
# [pharyngealized-long, long, short, pharyngealized] grouping per vowel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

function Set-FeatureRow {
    param($rowNum, $values)
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range($cols[$i] + $rowNum).Value = $values[$i]
    }
}

# --- Insert "iʶː" above row 35 ("iː") ---
$ws.Rows(35).Insert()
Set-FeatureRow 35 @("iʶː","+","-","+","+","-","-","+","-",0,"-","-","+",0,"-",0,"-","-",0,"-","-","-","+","+","-","-","+","+","-","+")

# --- Insert "uʶː" above the "uː" row (now row 39 after the insert above) ---
$ws.Rows(39).Insert()
Set-FeatureRow 39 @("uʶː","+","-","+","+","-","-","+","-",0,"-","-","+",0,"-","+","+","-",0,"-","-","-","+","+","-","+","+","+","-","+")

# --- Insert "aʶː" above the "aː" row (now row 43 after the two inserts above) ---
$ws.Rows(43).Insert()
Set-FeatureRow 43 @("aʶː","+","-","+","+","-","-","+","-",0,"-","-","+",0,"-",0,"-","-",0,"-","-","-","+","-","+","-","-","+","-","+")

# The shifted "iː" row (now row 36) also had its atr/rtr flags updated.
$ws.Range("AC36").Value = "-"
$ws.Range("AD36").Value = "+"

# Restore the view state (frozen pane position + active selection) to match
# where the author had scrolled/selected when saving.
$ws.Application.ActiveWindow.ScrollRow = 20
$sel = $ws.Range("AB35:AD35")
$sel.Select()
